# Updates the cryptocurrency price/volume table to the latest scrape.
# Generated to mirror the authoritative XML diff cell-by-cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as genuine text, even when it is a "numeric-looking"
# string (e.g. "1.002", "0.4827"). Excel COM auto-converts bare numeric
# strings to real numbers on assignment, which would both change the cell
# type away from text and introduce floating-point noise (1.002 ->
# 1.0019999999...). Marking the cell as Text first prevents that, and
# resetting the style back to Normal afterwards keeps the cell from
# picking up a stray number-format style the source file never had.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.941.43"
$ws.Range("E2").Value = "  -3.64%  "
Set-TextValue $ws.Range("D3") "1.722.88"
$ws.Range("E3").Value = "  -2.41%  "
Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  +0.11%  "
Set-TextValue $ws.Range("D5") "309.81"
$ws.Range("E5").Value = "  -5.67%  "
$ws.Range("E6").Value = "  +0.17%  "
Set-TextValue $ws.Range("D7") "0.4827"
$ws.Range("E7").Value = "  +3.77%  "
Set-TextValue $ws.Range("D8") "0.3465"
$ws.Range("E8").Value = "  -1.35%  "
Set-TextValue $ws.Range("D9") "43.17"
$ws.Range("E9").Value = "  -0.38%  "
Set-TextValue $ws.Range("D10") "0.07230"
$ws.Range("E10").Value = "  -1.83%  "
Set-TextValue $ws.Range("D11") "1.047"
$ws.Range("E11").Value = "  -3.19%  "
Set-TextValue $ws.Range("D12") "1.002"
$ws.Range("E12").Value = "  +0.16%  "
Set-TextValue $ws.Range("D13") "19.80"
$ws.Range("E13").Value = "  -4.05%  "
$ws.Range("E14").Value = "  -2.30%  "
Set-TextValue $ws.Range("D15") "1.716.64"
$ws.Range("E15").Value = "  -2.60%  "
Set-TextValue $ws.Range("D16") "6.778"
$ws.Range("E16").Value = "  -5.32%  "
Set-TextValue $ws.Range("D17") "86.76"
$ws.Range("E17").Value = "  -6.14%  "
Set-TextValue $ws.Range("D18") "0.00001033"
$ws.Range("E18").Value = "  -1.92%  "
Set-TextValue $ws.Range("D19") "0.06392"
$ws.Range("E19").Value = "  -0.50%  "
Set-TextValue $ws.Range("D20") "1.002"
$ws.Range("E21").Value = "  -2.01%  "
Set-TextValue $ws.Range("D22") "5.687"
$ws.Range("E22").Value = "  -1.38%  "
Set-TextValue $ws.Range("D23") "27.017.07"
$ws.Range("E23").Value = "  -3.46%  "
Set-TextValue $ws.Range("D24") "10.92"
$ws.Range("E24").Value = "  -1.98%  "
Set-TextValue $ws.Range("D25") "2.066"
$ws.Range("E25").Value = "  -4.15%  "
Set-TextValue $ws.Range("D26") "154.18"
$ws.Range("E26").Value = "  -4.98%  "
Set-TextValue $ws.Range("D27") "19.83"
$ws.Range("E27").Value = "  -1.04%  "
Set-TextValue $ws.Range("D28") "1.927.24"
$ws.Range("E28").Value = "  -1.89%  "
Set-TextValue $ws.Range("D29") "2.063"
$ws.Range("E29").Value = "  -4.88%  "
Set-TextValue $ws.Range("D30") "120.68"
$ws.Range("E30").Value = "  -1.79%  "
Set-TextValue $ws.Range("D31") "1.028"
$ws.Range("E31").Value = "  -3.91%  "
Set-TextValue $ws.Range("D32") "0.09286"
$ws.Range("E32").Value = "  +0.10%  "
Set-TextValue $ws.Range("D33") "3.612"
$ws.Range("E33").Value = "  -0.89%  "
Set-TextValue $ws.Range("D34") "5.359"
$ws.Range("E34").Value = "  -3.60%  "
Set-TextValue $ws.Range("D35") "0.05934"
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("E36").Value = "  -4.41%  "
Set-TextValue $ws.Range("D37") "1.428"
$ws.Range("E37").Value = "  +3.78%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D38") "10.89"
$ws.Range("E38").Value = "  -6.81%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D39") "0.1985"
$ws.Range("E39").Value = "  -3.64%  "
$ws.Range("E40").Value = "  +0.16%  "
Set-TextValue $ws.Range("D41") "4.719"
$ws.Range("E41").Value = "  -3.98%  "
Set-TextValue $ws.Range("D42") "0.5959"
$ws.Range("E42").Value = "  -2.78%  "
Set-TextValue $ws.Range("D43") "1.113"
$ws.Range("E43").Value = "  -5.74%  "
Set-TextValue $ws.Range("D44") "7.437"
$ws.Range("E44").Value = "  -4.43%  "
Set-TextValue $ws.Range("D45") "12.81"
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("E46").Value = "  -4.31%  "
$ws.Range("E47").Value = "  -3.59%  "
Set-TextValue $ws.Range("D48") "118.98"
$ws.Range("E48").Value = "  -3.20%  "
Set-TextValue $ws.Range("D49") "1.840"
$ws.Range("E49").Value = "  -4.44%  "
$ws.Range("E50").Value = "  -2.14%  "
Set-TextValue $ws.Range("D51") "0.06633"
$ws.Range("E51").Value = "  -2.67%  "
